$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, shifting the former row 32
# (Provincia de Limarí / 2021-12-22 entry) down to row 33.
$ws.Rows(32).Insert()

# New row 32 gets the data that row 31 used to hold, before row 31
# itself is updated to the new weekly entry below.
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44349
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112030
$ws.Range("G32").Value = "Poroto granado"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 60
$ws.Range("K32").Value = 30000
$ws.Range("L32").Value = 32000
$ws.Range("M32").Value = 31000
$ws.Range("N32").Value = "`$/saco 25 kilos"
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 1240
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"

# Row 31 is updated in place with the new weekly price entry.
$ws.Range("D31").Value = 44615
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 22000
$ws.Range("L31").Value = 23000
$ws.Range("M31").Value = 22500
$ws.Range("O31").Value = "Región del Maule"
$ws.Range("P31").Value = 900
